# Bump the published version 1.8.1 -> 1.8.2 on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.2"

# Add a new "Include from Tipo Identificad" worksheet after the existing
# "Include from identifierType" sheet. It reuses the same layout (Codes /
# All codes / System URI rows) as the other "Include from ..." sheets, so
# clone one of them to pick up identical formatting/styles/column widths
# and then just point its System URI cell at the new CodeSystem.
$src = $wb.Worksheets.Item("Include from identifierType")
$src.Copy($null, $src)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "Include from Tipo Identificad"
$newSheet.Range("B4").Value = "https://hl7chile.cl/fhir/ig/clcore/CodeSystem/CSTipoIdentificador"
